$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "ItemName"
$ws.Range("D1").Value = "Qty1"
$ws.Range("E1").Value = "ItemName2"
$ws.Range("F1").Value = "Qty2"
$ws.Range("G1").Value = "ItemName3"
$ws.Range("H1").Value = "Qty3"
$ws.Range("I1").Value = "Vendor"

# Row 2
$ws.Range("A2").Value = "Burger"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = ""
$ws.Range("I2").Value = "Uber eats"

# Row 3
$ws.Range("A3").Value = "Car"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = ""
$ws.Range("I3").Value = "Carvana"

# Row 4
$ws.Range("A4").Value = "Flowers"
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = ""
$ws.Range("I4").Value = "Target"

# Row 5
$ws.Range("A5").Value = "Watch"
$ws.Range("B5").Value = "Olivia"
$ws.Range("C5").Value = "James"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = ""
$ws.Range("I5").Value = "Amazon"

# Delete row 6 entirely
$ws.Rows("6").Delete()

# Column widths
$ws.Columns("A").ColumnWidth = 9.06640625
$ws.Columns("E").ColumnWidth = 11.9296875
$ws.Columns("G").ColumnWidth = 11.3984375

# Selection
$ws.Range("A6").Select()
